$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the previously-blank Start/End Time values (0) for a handful
#     of existing rows that only had a Date recorded so far. ---
$ws.Range("B60").Value = 0
$ws.Range("C60").Value = 0

$ws.Range("B72").Value = 0
$ws.Range("C72").Value = 0

$ws.Range("B76").Value = 0
$ws.Range("C76").Value = 0

$ws.Range("B77").Value = 0
$ws.Range("C77").Value = 0

$ws.Range("B78").Value = 0
$ws.Range("C78").Value = 0

# --- Add a new daily power record row (row 79) ---
$ws.Range("A79").Value = 43403
$ws.Range("D79").Formula = "=(C79-B79)* 1440"
$ws.Range("E79").Formula = "=IF(C79>B79, (C79-B79)*1440, (B79-C79)*1440)"
$ws.Range("F79").Formula = "=ABS((C79-B79)*1440)"

# --- Grow the table (comforter_cda_table) so the new row is included ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F79"))

# --- Match the saved selection/view state from the source workbook ---
$ws.Range("A79").Select()
